# Updates cryptos list D (Price) / E (Volume(1h)) columns for rows 2-51
# per the latest GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.279.66"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").Value = "1.677.02"
$ws.Range("E3").Value = "  +0.75%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.33"
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5285"
$ws.Range("E6").Value = "  +3.95%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.008"
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2687"
$ws.Range("E8").Value = "  +2.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06467"
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.89"
$ws.Range("E10").Value = "  +0.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07511"
$ws.Range("E11").Value = "  +1.47%  "
$ws.Range("D12").Value = "1.707.79"
$ws.Range("E12").Value = "  +2.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.513"
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5771"
$ws.Range("E14").Value = "  -0.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.000008499"
$ws.Range("E15").Value = "  +0.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.63"
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("D17").Value = "26.307.93"
$ws.Range("E17").Value = "  +1.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.911"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "189.86"
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.182"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.009"
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "144.80"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1273"
$ws.Range("E25").Value = "  +7.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.800"
$ws.Range("E26").Value = "  +2.84%  "
$ws.Range("E27").Value = "  +1.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.06497"
$ws.Range("E28").Value = "  -1.29%  "
$ws.Range("E29").Value = "  +3.93%  "
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("E31").Value = "  +1.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.586"
$ws.Range("E32").Value = "  +2.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.654"
$ws.Range("E33").Value = "  +1.57%  "
$ws.Range("E34").Value = "  +1.60%  "
$ws.Range("E35").Value = "  +2.27%  "
$ws.Range("E36").Value = "  +1.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.730"
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("E38").Value = "  +1.38%  "
$ws.Range("D39").Value = "1.114.63"
$ws.Range("E39").Value = "  +4.07%  "
$ws.Range("E40").Value = "  +1.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8719"
$ws.Range("E41").Value = "  +1.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.015"
$ws.Range("E42").Value = "  +0.71%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.44"
$ws.Range("E43").Value = "  +0.19%  "
$ws.Range("D44").Value = "1.827.29"
$ws.Range("E44").Value = "  +0.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00000000111"
$ws.Range("E45").Value = "  -2.21%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.92"
$ws.Range("E46").Value = "  +1.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.008"
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.163"
$ws.Range("E48").Value = "  +1.77%  "
$ws.Range("E49").Value = "  +1.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4292"
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.053"
$ws.Range("E51").Value = "  +2.03%  "
